$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 2020
$ws.Range("L4").Value = 10
$ws.Range("L5").Value = 6.4
$ws.Range("L6").Value = 13.5
$ws.Range("L7").Value = 24.3
$ws.Range("L8").Value = 27.8
$ws.Range("L9").Value = 20.9
$ws.Range("L10").Value = 26.7
$ws.Range("L11").Value = 28.4
$ws.Range("L12").Value = 25
